$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (column D) values for various coin rows as per the
# "Updated symbol list" commit.
$priceUpdates = @{
    "D2" = "268.09"
    "D3" = "22.88"
    "D4" = "6.306"
    "D5" = "0.06177"
    "D6" = "3.585"
    "D7" = "6.703"
    "D8" = "1.380"
    "D9" = "0.8414"
    "D10" = "0.01365"
    "D11" = "0.1603"
    "D12" = "0.08237"
    "D13" = "0.03412"
    "D14" = "0.03260"
    "D15" = "0.09280"
    "D16" = "3.910"
    "D17" = "0.001718"
    "D18" = "0.04863"
    "D19" = "0.006219"
    "D20" = "0.005363"
    "D23" = "3.759"
    "D24" = "2.326"
    "D25" = "0.3347"
    "D26" = "0.1212"
    "D40" = "0.04673"
    "D41" = "0.006963"
    "D42" = "0.1151"
    "D43" = "0.003602"
    "D44" = "0.01219"
    "D45" = "0.00006229"
    "D46" = "0.00000000750"
    "D47" = "0.7004"
    "D48" = "0.1766"
    "D49" = "0.00001401"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so values remain exact strings (e.g. "268.09")
    # rather than being coerced into numeric cells.
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    # Restore default style so no stray number format is left behind.
    $cell.Style = "Normal"
}
